# Lec 8/9 slide deck refresh:
#  - the "automatic date" placeholders on the slide master, every slide
#    layout, and the notes master were refreshed by PowerPoint to the
#    new save date (6.05.2023 -> 10.05.2024)
#  - a small wording fix on the "Endeksler" slide
#    ("otomatik endekslerler" -> "otomatik endeksler")

$p = $ppt.ActivePresentation
$oldDate = "6.05.2023"
$newDate = "10.05.2024"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes, $label) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $tr = $sh.TextFrame.TextRange
                if ($tr.Text -eq $oldDate) {
                    $tr.Text = $newDate
                }
            }
        }
    }
}

# 1) Slide master date placeholder
$masterLabel = "Master"
Update-DatePlaceholder $p.SlideMaster.Shapes $masterLabel

# 2) Every slide layout's date placeholder
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $cl = $layouts.Item($li)
    $layoutLabel = "Layout" + $li
    Update-DatePlaceholder $cl.Shapes $layoutLabel
}

# 3) Notes master date placeholder
$notesLabel = "NotesMaster"
Update-DatePlaceholder $p.NotesMaster.Shapes $notesLabel

# 4) Wording fix on slide 5 ("Endeksler"): "otomatik endekslerler" -> "otomatik endeksler"
$slide = $p.Slides.Item(5)
$shape = $slide.Shapes.Item(2)
$tr = $shape.TextFrame.TextRange
$oldPhrase = " nitelikleri otomatik endekslerler."
$newPhrase = " nitelikleri otomatik endeksler."
$para = $tr.Paragraphs(2, 1)
$localIdx = $para.Text.IndexOf($oldPhrase)
if ($localIdx -ge 0) {
    $absStart = $para.Start + $localIdx
    $run = $tr.Characters($absStart, $oldPhrase.Length)
    $run.Text = $newPhrase
}
